$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.248.33"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "2.255.66"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "80.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.22%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "2.585.07"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "2.228.58"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").Value = "44.128.40"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("E22").Value = "  +9.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.47%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0880"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.86%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0370"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.16%  "
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.37%  "
$ws.Range("E40").Value = "  +21.21%  "
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.18%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.59%  "
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.460"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +27.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.73%  "
